$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns/cells that are no longer part of the table (C and D)
$ws.Range("C1:D2").Clear()

# Update header text
$ws.Range("B1").Value = "Analysis"

# Row 2 - update B2's value (A2 stays 0, already styled)
$ws.Range("B2").Value = "{'Doopa': 'dopa'}"

# Row 3 - new row, mirror A2's style via copy/paste-special (formats only)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "{'Doopa': 'dopa'}"

# Row 4 - new row
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "      HER_#1_#1.DTA  HER_#1_#2.DTA  HER_#1_#3.DTA  HER_#2_#1.DTA  HER_#2_#2.DTA  ...  HER_POWROT_#2_#1.DTA  HER_POWROT_#2_#2.DTA  HER_POWROT_#2_#3.DTA      mean       std`n0.01      -0.015307       -0.01532      -0.015309      -0.015348      -0.015335  ...               -0.0153             -0.015266             -0.015319 -0.015302  0.000027`n[1 rows x 14 columns]"

$excel.CutCopyMode = $false
